# Herbie usability-tests sheet: reorder the "Webchart Prescribe a Medicine"
# (id 3) and "ADP Task" (id 4) rows so the ADP task comes first, and drop
# the stray formatted-but-empty row that trailed the data.
#
# Plan:
#   1. Swap the A4:E5 <-> A5:E5 cell content (values + per-cell formatting)
#      through a staging range, since a straight row-swap isn't available.
#   2. Row height lives on the row, not the cell range, so swap it too.
#   3. Delete the now-redundant empty row 6.
#   4. Hyperlinks are anchored to a fixed cell and don't follow the cells
#      they describe when copied, so rebuild the four hyperlinks pointing
#      at their new homes (same URLs/anchor order as the finished sheet).
#   5. Adding a hyperlink re-stamps the cell's font, so re-pin each D-cell
#      back onto the workbook's "Hyperlink" style (plus the original
#      vertical-top alignment on D5) to match the pre-existing formatting.
#   6. Leave the selection on E5, matching where the edit ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap the row 4 / row 5 contents via a far-away staging range.
$null = $ws.Range("A4:E4").Copy($ws.Range("A100:E100"))
$null = $ws.Range("A5:E5").Copy($ws.Range("A4:E4"))
$null = $ws.Range("A100:E100").Copy($ws.Range("A5:E5"))
$null = $ws.Range("A100:E100").Clear()

# 2) Row height is row-level, so swap it explicitly.
$origHeight4 = $ws.Rows(4).RowHeight()
$origHeight5 = $ws.Rows(5).RowHeight()
$ws.Rows(4).RowHeight = $origHeight5
$ws.Rows(5).RowHeight = $origHeight4

# 3) Drop the leftover empty row.
$null = $ws.Rows(6).Delete()

# 4) Rebuild the hyperlinks against their new anchor cells.
$null = $ws.Range("A1").Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("D5"), "https://masterdaily.dev.webchart.app/webchart.cgi")
$null = $ws.Hyperlinks.Add($ws.Range("D2"), "http://18.219.226.248/academy.html")
$null = $ws.Hyperlinks.Add($ws.Range("D4"), "https://workforcenow.adp.com/theme/index.html", "/home")
$null = $ws.Hyperlinks.Add($ws.Range("D3"), "https://masterdaily.dev.webchart.app/webchart.cgi")

# 5) Restore the cell styling that Hyperlinks.Add overwrote.
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("D3").Style = "Hyperlink"
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("D5").Style = "Hyperlink"
$ws.Range("D5").VerticalAlignment = -4160

# 6) Match the active selection left behind by the edit.
$null = $ws.Range("E5").Select()
